$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "65÷8=8, 1"
$t.Cell(1, 2).Range.Text = "81÷2=40, 1"
$t.Cell(1, 3).Range.Text = "78÷3=26, 0"
$t.Cell(1, 4).Range.Text = "52÷5=10, 2"
$t.Cell(1, 5).Range.Text = "54÷6=9, 0"

$t.Cell(5, 1).Range.Text = "81÷3=27, 0"
$t.Cell(5, 2).Range.Text = "59÷7=8, 3"
$t.Cell(5, 3).Range.Text = "81÷2=40, 1"
$t.Cell(5, 4).Range.Text = "87÷5=17, 2"
$t.Cell(5, 5).Range.Text = "17÷3=5, 2"

$t.Cell(9, 1).Range.Text = "91÷7=13, 0"
$t.Cell(9, 2).Range.Text = "96÷4=24, 0"
$t.Cell(9, 3).Range.Text = "20÷2=10, 0"
$t.Cell(9, 4).Range.Text = "75÷4=18, 3"
$t.Cell(9, 5).Range.Text = "59÷9=6, 5"

$t.Cell(13, 1).Range.Text = "49÷4=12, 1"
$t.Cell(13, 2).Range.Text = "89÷4=22, 1"
$t.Cell(13, 3).Range.Text = "12÷8=1, 4"
$t.Cell(13, 4).Range.Text = "70÷2=35, 0"
$t.Cell(13, 5).Range.Text = "47÷9=5, 2"

$t.Cell(17, 1).Range.Text = "80÷9=8, 8"
$t.Cell(17, 2).Range.Text = "17÷3=5, 2"
$t.Cell(17, 3).Range.Text = "34÷3=11, 1"
$t.Cell(17, 4).Range.Text = "48÷3=16, 0"
$t.Cell(17, 5).Range.Text = "22÷4=5, 2"
